# "Generate Report for Handback" — localization-status.xlsx update
#
# Summary of the change being applied:
#  - The "Ready for handoff" status (used on the Overview sheet's zh-cn/de-de
#    columns and on each language sheet's Status cell) is now
#    "Handed back: in sync with en-US" everywhere it appears.
#  - The zh-cn and de-de sheets now have their "Latest Target File" (I2) and
#    "Latest Handback File" (J2) populated (I2 becomes a hyperlink to the
#    source markdown file, J2 holds the generated xliff file name), and the
#    "Latest Handback DateTime" (K2) is stamped with the handback time.
#  - A handful of columns are widened to comfortably fit the newly
#    populated long file-name / status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"
$mdFileName = "e1a71351-3ca9-4e52-a4cf-3d9aa822e7a4.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6decf80307c617017b1074c44e636dab1e40da5b/e2e/e1a71351-3ca9-4e52-a4cf-3d9aa822e7a4.md"
$zhXlf      = "e1a71351-3ca9-4e52-a4cf-3d9aa822e7a4.cebe26dfc03180e89d7c64921921cd2b181e439b.zh-cn.xlf"
$deXlf      = "e1a71351-3ca9-4e52-a4cf-3d9aa822e7a4.cebe26dfc03180e89d7c64921921cd2b181e439b.de-de.xlf"

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F2 and the Status cell on each language sheet, C2)
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$zhcn.Range("C2").Value = $handedBack
$dede.Range("C2").Value = $handedBack

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (I2), Latest Handback File (J2)
#    and Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$zhcn.Range("J2").Value = $zhXlf
$zhcn.Range("K2").Value = "2016-08-24 09:01:56"

$zhI2 = $zhcn.Range("I2")
$zhcn.Hyperlinks.Add($zhI2, $mdUrl, "", "", $mdFileName)
$zhI2.Style = "HyperLink"

# ---------------------------------------------------------------------------
# 3. de-de sheet: fill in Latest Target File (I2), Latest Handback File (J2)
#    and Latest Handback DateTime (K2)
# ---------------------------------------------------------------------------
$dede.Range("J2").Value = $deXlf
$dede.Range("K2").Value = "2016-08-24 09:02:13"

$deI2 = $dede.Range("I2")
$dede.Hyperlinks.Add($deI2, $mdUrl, "", "", $mdFileName)
$deI2.Style = "HyperLink"

# ---------------------------------------------------------------------------
# 4. Widen columns that now show the longer status text / file names
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.15   # E -> ~29.98 stored width
$overview.Columns.Item(6).ColumnWidth = 29.15   # F -> ~29.98 stored width

$zhcn.Columns.Item(3).ColumnWidth = 29.15        # C -> ~29.98 stored width
$zhcn.Columns.Item(9).ColumnWidth = 39.1         # I -> 40
$zhcn.Columns.Item(10).ColumnWidth = 39.1        # J -> 40

$dede.Columns.Item(3).ColumnWidth = 29.15        # C -> ~29.98 stored width
$dede.Columns.Item(9).ColumnWidth = 39.1         # I -> 40
$dede.Columns.Item(10).ColumnWidth = 39.1        # J -> 40
